$d = $word.ActiveDocument

# ---------------------------------------------------------------
# "Anime Buzz" project entry:
#   "Anime Buzz"
#     -> "Anime Buzz - " (unchanged formatting)
#        + "github.com/jakemoritz/Anime-Buzz" (new, bold run)
# ---------------------------------------------------------------

$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Anime Buzz", $true, $true, $false, $false, $false, $true, 1, $false, `
    "Anime Buzz - github.com/jakemoritz/Anime-Buzz", 2)

$rng1b = $d.Content
$found1b = $rng1b.Find.Execute("github.com/jakemoritz/Anime-Buzz", $true, $true, $false, $false, $false, $true, 1, $false, `
    "", 0)
if ($found1b) {
    $rng1b.Font.Bold = $true
}

# ---------------------------------------------------------------
# "Tasking" project entry:
#   "Tasking  - github.com/jakemoritz/Tasking"
#     -> "Tasking  -" (unchanged formatting)
#        + " github.com/jakemoritz/Tasking" (new, bold run)
# ---------------------------------------------------------------

$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Tasking  - github.com/jakemoritz/Tasking", $true, $true, $false, $false, $false, $true, 1, $false, `
    "Tasking  - github.com/jakemoritz/Tasking", 2)

$rng2b = $d.Content
$found2b = $rng2b.Find.Execute(" github.com/jakemoritz/Tasking", $true, $true, $false, $false, $false, $true, 1, $false, `
    "", 0)
if ($found2b) {
    $rng2b.Font.Bold = $true
}

Write-Host "AnimeBuzz split: $found1 / bolded: $found1b"
Write-Host "Tasking split: $found2 / bolded: $found2b"
